# "Add files via upload" — re-upload of geiger.xlsx with the "Element"
# sheet's B89:B105 column switched on (0 -> 1) and the sheet scrolled /
# selected down around that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Element")

# Make sure we're working on the right sheet.
$ws.Activate()

# B89:B105 changed from 0 to 1.
$ws.Range("B89:B105").Value = 1

# Reflect the new selection / scroll position recorded in the saved view:
# top-left visible cell around row 70 and the active selection on B89:B105.
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B89:B105").Select()
